# Aggiornamento a l 23 agosto 2021
# Appends the new daily rows (2021-08-10 .. 2021-08-23) to the bottom of the
# single data table on the active sheet, extending it from row 343 to row 357.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastExistingRow = 343

# Copy the formatting (date style, borders, bold, alignment, etc.) of the
# last existing data row down across the 14 new rows before writing values,
# so column A keeps its date/time number format + style (same as every
# other row in the table).
$ws.Range("A343:D343").Copy()
$ws.Range("A344:D357").PasteSpecial(-4122)

# r, date-serial (col A), nuovi pos. (col B), somma mobile 7gg. (col C),
# somma mobile 7gg. per 100mila abitanti (col D)
$data = @(
    @(344, 44418, 0, 2, 57.75339301183945),
    @(345, 44419, 0, 2, 57.75339301183945),
    @(346, 44420, 1, 3, 86.63008951775916),
    @(347, 44421, 0, 2, 57.75339301183945),
    @(348, 44422, 0, 2, 57.75339301183945),
    @(349, 44423, 2, 3, 86.63008951775916),
    @(350, 44424, 0, 3, 86.63008951775916),
    @(351, 44425, 0, 3, 86.63008951775916),
    @(352, 44426, 0, 3, 86.63008951775916),
    @(353, 44427, 1, 3, 86.63008951775916),
    @(354, 44428, 5, 8, 231.0135720473578),
    @(355, 44429, 0, 8, 231.0135720473578),
    @(356, 44430, 1, 7, 202.1368755414381),
    @(357, 44431, 0, 7, 202.1368755414381)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
